# Initial DBF to CSV export
# Populate "variables" sheet (NIJ_CFS_VARS bank) with the new calls-for-service fields.
$wb = $excel.ActiveWorkbook

$wsVars = $wb.Worksheets.Item("variables")
$wsRecords = $wb.Worksheets.Item("records")
$wsLayouts = $wb.Worksheets.Item("layouts")

# --- variables sheet: names first (column C), then datatypes (column E) ---
$wsVars.Range("C2").Value = "category"
$wsVars.Range("C3").Value = "call_group"
$wsVars.Range("C4").Value = "final_case"
$wsVars.Range("C5").Value = "case_desc"
$wsVars.Range("C6").Value = "occ_date"

$wsVars.Range("E2").Value = "text"
$wsVars.Range("E3").Value = "text"
$wsVars.Range("E4").Value = "text"
$wsVars.Range("E5").Value = "text"

$wsVars.Range("E7").Value = "numeric"
$wsVars.Range("E8").Value = "numeric"
$wsVars.Range("E9").Value = "numeric"

# --- records sheet: register the new cfs_dbf record ---
$wsRecords.Range("A2").Value = "NIJ_CFS_RL"
$wsRecords.Range("A2").VerticalAlignment = -4160

# --- variables sheet: bank column ---
$wsVars.Range("A2").Value = "NIJ_CFS_VARS"
$wsVars.Range("A3").Value = "NIJ_CFS_VARS"
$wsVars.Range("A4").Value = "NIJ_CFS_VARS"
$wsVars.Range("A5").Value = "NIJ_CFS_VARS"
$wsVars.Range("A6").Value = "NIJ_CFS_VARS"
$wsVars.Range("A7").Value = "NIJ_CFS_VARS"
$wsVars.Range("A8").Value = "NIJ_CFS_VARS"
$wsVars.Range("A9").Value = "NIJ_CFS_VARS"

$wsRecords.Range("E2").Value = "NIJ_CFS_VARS"
$wsRecords.Range("E2").VerticalAlignment = -4160
$wsRecords.Range("B2").Value = "cfs_dbf"
$wsRecords.Range("C2").Value = "calls-for-service DBF"

# --- variables sheet: remaining datatype + names ---
$wsVars.Range("E6").Value = "date(iso8601)"
$wsVars.Range("C7").Value = "x_coordina"
$wsVars.Range("C8").Value = "y_coordina"
$wsVars.Range("C9").Value = "census_tra"

$wsVars.Range("F7").Value = 0
$wsVars.Range("F8").Value = 0
$wsVars.Range("F9").Value = 0

# --- layouts sheet: cfs_dbf field layout rows ---
$wsLayouts.Range("A2:A9").Style = "Normal"
$wsLayouts.Range("A2").Value = "cfs_dbf"
$wsLayouts.Range("B2").Value = "category"
$wsLayouts.Range("A3").Value = "cfs_dbf"
$wsLayouts.Range("B3").Value = "call_group"
$wsLayouts.Range("A4").Value = "cfs_dbf"
$wsLayouts.Range("B4").Value = "final_case"
$wsLayouts.Range("A5").Value = "cfs_dbf"
$wsLayouts.Range("B5").Value = "case_desc"
$wsLayouts.Range("A6").Value = "cfs_dbf"
$wsLayouts.Range("B6").Value = "occ_date"
$wsLayouts.Range("A7").Value = "cfs_dbf"
$wsLayouts.Range("B7").Value = "x_coordina"
$wsLayouts.Range("A8").Value = "cfs_dbf"
$wsLayouts.Range("B8").Value = "y_coordina"
$wsLayouts.Range("A9").Value = "cfs_dbf"
$wsLayouts.Range("B9").Value = "census_tra"

# --- selections / active sheet / window ---
[void]$wsVars.Range("C2:C9").Select()
[void]$wsRecords.Range("B2").Select()
[void]$wsLayouts.Range("D14").Select()
[void]$wsLayouts.Activate()
